# CU Administrar Cliente y la interfaz alumno
# Se realizó el CU Administrar Cliente y la interfaz alumno

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# --- Update task status (column F) from "En proceso" to "Hecho" ---
$ws.Range("F8").Value  = "Hecho"
$ws.Range("F11").Value = "Hecho"
$ws.Range("F16").Value = "Hecho"
$ws.Range("F23").Value = "Hecho"

# --- Register daily consumed hours (fills in previously-blank "Cons." cells) ---
# Row 8 - "Realizar inferfaces en FW" (Alumno) - Día 14
$ws.Range("AU8").Value = 1

# Row 11 - "Crear el diagrama de clases" - Día 6, Día 8, Día 9
$ws.Range("W11").Value  = 2
$ws.Range("AC11").Value = 1
$ws.Range("AF11").Value = 2

# Row 16 - "Generar las pruebas de la interfaz Cliente" - Día 19
$ws.Range("BJ16").Value = 0.5

# Row 18 - "Generar las pruebas de la interfaz Profesor" - Día 19
$ws.Range("BJ18").Value = 0.5

# Row 23 - "Dar de baja alumno" - Día 19
$ws.Range("BJ23").Value = 2

# --- Update the frozen-pane view / selection to reflect the latest work area ---
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 7
$ws.Range("D25").Select()

$wb.Save()
